# Update "想去人数" (interested-count) figures to the latest scrape values.
# (sheet index, row, new value) — sheet order per workbook.xml:
#   1 = 展览, 2 = 演出, 3 = 本地生活, 4 = 全部类型

$wb = $excel.ActiveWorkbook

function Set-F {
    param($SheetIdx, $Row, $NewValue)
    $ws = $wb.Worksheets.Item($SheetIdx)
    $ws.Cells.Item($Row, 6).Value = $NewValue
}

# 展览 (sheet 1)
Set-F 1 6  2417
Set-F 1 8  1787
Set-F 1 9  3034
Set-F 1 11 4517
Set-F 1 12 400
Set-F 1 13 224
Set-F 1 15 570
Set-F 1 18 250
Set-F 1 20 115
Set-F 1 21 315
Set-F 1 22 4553
Set-F 1 24 4098
Set-F 1 27 595
Set-F 1 28 4388
Set-F 1 30 651
Set-F 1 31 619
Set-F 1 32 597

# 演出 (sheet 2)
Set-F 2 4 6

# 本地生活 (sheet 3)
Set-F 3 3 1047

# 全部类型 (sheet 4)
Set-F 4 4  1047
Set-F 4 9  2417
Set-F 4 11 1787
Set-F 4 13 3034
Set-F 4 15 4517
Set-F 4 16 400
Set-F 4 17 224
Set-F 4 19 570
Set-F 4 22 250
Set-F 4 25 115
Set-F 4 26 315
Set-F 4 27 4553
Set-F 4 29 4098
Set-F 4 32 595
Set-F 4 33 4388
Set-F 4 35 651
Set-F 4 36 619
Set-F 4 37 597
Set-F 4 38 6
